$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("Donation")
$ws4.Range("C2").Value = 123456
$ws4.Range("D2").Value = "Ateneo"
$ws4.Range("F2").Value = 654321
$ws4.Range("G2").Value = (Get-Date -Year 2021 -Month 1 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws4.Range("G2").NumberFormat = "m/d/yyyy"
$ws4.Range("H2").Value = 20000
$ws4.Range("I2").Value = "nothing"
$ws4.Range("M2").Value = "yes"
$ws4.Range("N2").Value = "test"
Write-Output "done"
